$d = $word.ActiveDocument

# 1. Date day: 24 -> 28 (only the first occurrence, the standalone day number;
#    "24" also appears inside "2024" further along so ReplaceOne is required)
$d.Content.Find.Execute("24", $true, $false, $false, $false, $false, $true, 1, $false, "28", 1)

# 2. Month: " de marzo de 2024" -> " de abril de 2024"
$d.Content.Find.Execute(" de marzo de 2024", $true, $false, $false, $false, $false, $true, 1, $false, " de abril de 2024", 2)

# 3. "nuestra alumna Señorita" -> "nuestro alumno Señor"
$d.Content.Find.Execute("nuestra alumna Señorita", $true, $false, $false, $false, $false, $true, 1, $false, "nuestro alumno Señor", 2)

# 4. RUT: 21061253-3 -> 21061253-K
$d.Content.Find.Execute("21061253-3", $true, $false, $false, $false, $false, $true, 1, $false, "21061253-K", 2)

# 5. ", en calidad de alumna" -> ", en calidad de alumno"
$d.Content.Find.Execute(", en calidad de alumna", $true, $false, $false, $false, $false, $true, 1, $false, ", en calidad de alumno", 2)

# 6. "La señorita" -> "El señor"
$d.Content.Find.Execute("La señorita", $true, $false, $false, $false, $false, $true, 1, $false, "El señor", 2)

# 7. "Sexto Semestre" -> "Primer Semestre" (case-sensitive match so the later
#    lowercase "sexto semestre" further in the paragraph is left untouched)
$d.Content.Find.Execute("Sexto Semestre", $true, $false, $false, $false, $false, $true, 1, $false, "Primer Semestre", 2)

# 8. "LA ALUMNA" -> "EL ALUMNO"
$d.Content.Find.Execute("LA ALUMNA", $true, $false, $false, $false, $false, $true, 1, $false, "EL ALUMNO", 2)
